$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "308.28"
$ws.Range("E2").Value = "0.88%"
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "38.55"
$ws.Range("E3").Value = "7.90%"
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.103"
$ws.Range("E4").Value = "1.19%"
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.08131"
$ws.Range("E5").Value = "1.08%"
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "1.966"
$ws.Range("E6").Value = "3.70%"
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "7.939"
$ws.Range("E7").Value = "2.10%"
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.9284"
$ws.Range("E8").Value = "0.82%"
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "0.1451"
$ws.Range("E9").Value = "14.94%"
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.1958"
$ws.Range("E10").Value = "2.41%"
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.09198"
$ws.Range("E11").Value = "1.35%"
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.03501"
$ws.Range("E12").Value = "1.03%"
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.09804"
$ws.Range("E13").Value = "-0.52%"
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001400"
$ws.Range("E14").Value = "-1.06%"
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.006086"
$ws.Range("E15").Value = "-2.47%"
$rng.Style = "Normal"

$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$ws.Range("E16").Value = "-4.51%"
$rng.Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "4.200"
$ws.Range("E17").Value = "1.36%"
$rng.Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "3.454"
$ws.Range("E18").Value = "3.44%"
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = "0.82%"
$rng.Style = "Normal"

$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$ws.Range("E20").Value = "-2.67%"
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "4.802"
$ws.Range("E21").Value = "-7.12%"
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "0.2453"
$ws.Range("E22").Value = "6.46%"
$rng.Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.04429"
$ws.Range("E23").Value = "0.03%"
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").Value = "-1.20%"
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.004847"
$ws.Range("E25").Value = "5.06%"
$rng.Style = "Normal"

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.02116"
$ws.Range("E39").Value = "8.68%"
$rng.Style = "Normal"

$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$ws.Range("E40").Value = "-4.58%"
$rng.Style = "Normal"

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.007464"
$ws.Range("E41").Value = "-1.86%"
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.01014"
$ws.Range("E42").Value = "-0.55%"
$rng.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.1364"
$ws.Range("E43").Value = "0.96%"
$rng.Style = "Normal"

$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.002142"
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.009411"
$ws.Range("E45").Value = "-2.50%"
$rng.Style = "Normal"

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00006220"
$ws.Range("E46").Value = "1.61%"
$rng.Style = "Normal"

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$ws.Range("E47").Value = "-0.03%"
$rng.Style = "Normal"

$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.003067"
$rng.Style = "Normal"

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = "0.001601"
$ws.Range("E49").Value = "-3.54%"
$rng.Style = "Normal"

$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
$rng.Style = "Normal"

$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$ws.Range("E51").Value = "-0.03%"
$rng.Style = "Normal"
